$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 11.80884686099532
